$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($row, $col, $text) {
    # Force the cell to be stored as literal text (avoids Excel's automatic
    # number/date recognition for values like "01/18/2025" or "26"), then
    # restore the plain "Normal" style so no stray number format sticks
    # around on the cell.
    $cell = $ws.Cells.Item($row, $col)
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $cell.Style = "Normal"
}

# Row 54: today's first new DAP (Data-At-a-Point) reading. Columns A and D
# become real numbers; B/C stay as text.
$ws.Cells.Item(54, 1).Value = 26
Set-TextCell 54 2 "Partly Cloudy"
Set-TextCell 54 3 "01/18/2025"
$ws.Cells.Item(54, 4).Value = 21

# Row 55: next DAP reading for today.
$ws.Cells.Item(55, 1).Value = 26
Set-TextCell 55 2 "Partly Cloudy"
Set-TextCell 55 3 "01/18/2025"
$ws.Cells.Item(55, 4).Value = 22

# Row 56: next DAP reading for today.
$ws.Cells.Item(56, 1).Value = 26
Set-TextCell 56 2 "Partly Cloudy"
Set-TextCell 56 3 "01/18/2025"
$ws.Cells.Item(56, 4).Value = 23

# Row 57: latest DAP reading for today - recorded as plain text, matching
# how the scraper appends its most recent (not-yet-finalized) entry.
Set-TextCell 57 1 "26"
Set-TextCell 57 2 "Partly Cloudy"
Set-TextCell 57 3 "01/18/2025"
Set-TextCell 57 4 "23"
